$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 26 (pushes existing rows 26-33 down to 27-34),
# inheriting formatting (row height / styles) from the row above.
$ws1.Rows.Item(26).Insert() | Out-Null
$ws1.Rows.Item(26).RowHeight = 48.75

# Populate the new row. Columns are written Name, Description, Value so the
# shared-string table picks up the same ordering as the authored workbook.
$ws1.Range("A26").Value = "CountryLookUpPath"
$ws1.Range("C26").Value = "File path for Excel doc containing valid countries for certificates to be sent to. Stored in the shared drive. "
$ws1.Range("B26").Value = "\\EARTH.GSI.GOV.UK\USER\SHARED\Agency\CoFS for G drive\RobotDocuments\ValidCountriesList.xlsx"

# Grow Table1 so the new row is included in the table range.
$tbl = $ws1.ListObjects.Item(1)
$tbl.Resize($ws1.Range("A1:C34")) | Out-Null

# Journey testing: make Sheet1 the active sheet/tab and leave the selection on F17.
$ws1.Activate() | Out-Null
$ws1.Range("F17").Select() | Out-Null
